$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 249
$ws.Range("A249").Value = '2024-10-01 21:10:52'
$ws.Range("B249").Value = 'check_availability'
$ws.Range("C249").Value = 'https://example.com'
$ws.Range("D249").Value = 'Checked availability: Selected or default date current date is available for booking.'
$ws.Range("E249").Value = '2024-10-01'
$ws.Range("E249").NumberFormat = "@"
$ws.Range("E249").Value = '2024-10-01'
$ws.Range("E249").ClearFormats()
$ws.Range("F249").Value = '21:10:52'

# Row 250
$ws.Range("A250").Value = '2024-10-01 21:10:52'
$ws.Range("B250").Value = 'check_availability'
$ws.Range("C250").Value = 'https://example.com'
$ws.Range("D250").Value = 'Failed to check availability: Failed to check availability'
$ws.Range("E250").Value = '2024-10-01'
$ws.Range("E250").NumberFormat = "@"
$ws.Range("E250").Value = '2024-10-01'
$ws.Range("E250").ClearFormats()
$ws.Range("F250").Value = '21:10:52'

# Row 251
$ws.Range("A251").Value = '2024-10-01 21:10:52'
$ws.Range("B251").Value = 'check_availability'
$ws.Range("C251").Value = 'https://example.com'
$ws.Range("D251").Value = 'Checked availability: No availability for the selected date.'
$ws.Range("E251").Value = '2024-10-01'
$ws.Range("E251").NumberFormat = "@"
$ws.Range("E251").Value = '2024-10-01'
$ws.Range("E251").ClearFormats()
$ws.Range("F251").Value = '21:10:52'

# Row 252
$ws.Range("A252").Value = '2024-10-01 21:10:52'
$ws.Range("B252").Value = 'check_availability'
$ws.Range("C252").Value = 'https://example.com/product'
$ws.Range("D252").Value = '$199.99'
$ws.Range("D252").NumberFormat = "@"
$ws.Range("D252").Value = '$199.99'
$ws.Range("D252").ClearFormats()
$ws.Range("E252").Value = '2024-10-01'
$ws.Range("E252").NumberFormat = "@"
$ws.Range("E252").Value = '2024-10-01'
$ws.Range("E252").ClearFormats()
$ws.Range("F252").Value = '21:10:52'

# Row 253
$ws.Range("A253").Value = '2024-10-01 21:10:53'
$ws.Range("B253").Value = 'check_availability'
$ws.Range("C253").Value = 'invalid_url'
$ws.Range("D253").Value = 'Error fetching price: Invalid URL'
$ws.Range("E253").Value = '2024-10-01'
$ws.Range("E253").NumberFormat = "@"
$ws.Range("E253").Value = '2024-10-01'
$ws.Range("E253").ClearFormats()
$ws.Range("F253").Value = '21:10:53'

# Row 254
$ws.Range("A254").Value = '2024-10-01 21:10:53'
$ws.Range("B254").Value = 'check_availability'
$ws.Range("C254").Value = 'https://example.com'
$ws.Range("D254").Value = 'Checked availability: Selected or default date is available for booking.'
$ws.Range("E254").Value = '2024-10-01'
$ws.Range("E254").NumberFormat = "@"
$ws.Range("E254").Value = '2024-10-01'
$ws.Range("E254").ClearFormats()
$ws.Range("F254").Value = '21:10:53'

# Row 255
$ws.Range("A255").Value = '2024-10-01 21:10:54'
$ws.Range("B255").Value = 'check_availability'
$ws.Range("C255").Value = 'https://example.com'
$ws.Range("D255").Value = 'Failed to check availability: Failed to check availability'
$ws.Range("E255").Value = '2024-10-01'
$ws.Range("E255").NumberFormat = "@"
$ws.Range("E255").Value = '2024-10-01'
$ws.Range("E255").ClearFormats()
$ws.Range("F255").Value = '21:10:54'

# Row 256
$ws.Range("A256").Value = '2024-10-01 21:10:55'
$ws.Range("B256").Value = 'check_availability'
$ws.Range("C256").Value = 'https://example.com/product'
$ws.Range("D256").Value = '100 USD'
$ws.Range("E256").Value = '2024-10-01'
$ws.Range("E256").NumberFormat = "@"
$ws.Range("E256").Value = '2024-10-01'
$ws.Range("E256").ClearFormats()
$ws.Range("F256").Value = '21:10:55'
